$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new log entry (2026/02/07, 1 o'clock) needs to be inserted right after
# the last 2026/02/06 entry (row 782), before the 2026/12/29 block that
# currently starts at row 783. Insert a new row there, which pushes the
# old rows 783:824 down to 784:825 and grows the sheet to A1:D825.
$ws.Rows.Item(783).Insert()

# Force column A to be stored as plain text so the date-like string isn't
# auto-converted into a date serial value, then restore the default
# "Normal" style so the cell doesn't keep a lingering custom number format.
$ws.Cells.Item(783, 1).NumberFormat = "@"
$ws.Cells.Item(783, 1).Value = "2026/02/07"
$ws.Cells.Item(783, 1).Style = "Normal"

$ws.Cells.Item(783, 2).Value = "土"
$ws.Cells.Item(783, 3).Value = 1
$ws.Cells.Item(783, 4).Value = 201
